$wb = $excel.ActiveWorkbook

# --- Cover sheet: update "Latest Update" date cell C12 from text to a real date value ---
$cover = $wb.Worksheets.Item("Cover")
$cover.Range("C12").Value = 43872

# --- sa_all_tot sheet: refresh stock-addition (tonnes) figures in column F for rows 2-49 ---
$ws = $wb.Worksheets.Item("sa_all_tot")

$ws.Range("F2").Value = 100017445.1636125
$ws.Range("F2").Font.Name = "Arial"
$ws.Rows.Item(2).AutoFit()

$ws.Range("F3").Value = 106298643.6820022
$ws.Range("F3").Font.Name = "Arial"
$ws.Rows.Item(3).AutoFit()

$ws.Range("F4").Value = 106443874.80832499
$ws.Range("F4").Font.Name = "Arial"
$ws.Rows.Item(4).AutoFit()

$ws.Range("F5").Value = 38221683.543860734
$ws.Range("F5").Font.Name = "Arial"
$ws.Rows.Item(5).AutoFit()

$ws.Range("F6").Value = 553065039.60840642
$ws.Range("F6").Font.Name = "Arial"
$ws.Rows.Item(6).AutoFit()

$ws.Range("F7").Value = 392357947.85977751
$ws.Range("F7").Font.Name = "Arial"
$ws.Rows.Item(7).AutoFit()

$ws.Range("F8").Value = 64860705.785781913
$ws.Range("F8").Font.Name = "Arial"
$ws.Rows.Item(8).AutoFit()

$ws.Range("F9").Value = 13972337415.664881
$ws.Range("F9").Font.Name = "Arial"
$ws.Rows.Item(9).AutoFit()

$ws.Range("F10").Value = 11549774.61997688
$ws.Range("F10").Font.Name = "Arial"
$ws.Rows.Item(10).AutoFit()

$ws.Range("F11").Value = 91774596.56463705
$ws.Range("F11").Font.Name = "Arial"
$ws.Rows.Item(11).AutoFit()

$ws.Range("F12").Value = 700656399.86290944
$ws.Range("F12").Font.Name = "Arial"
$ws.Rows.Item(12).AutoFit()

$ws.Range("F13").Value = 52222075.784433663
$ws.Range("F13").Font.Name = "Arial"
$ws.Rows.Item(13).AutoFit()

$ws.Range("F14").Value = 12003488.616244361
$ws.Range("F14").Font.Name = "Arial"
$ws.Rows.Item(14).AutoFit()

$ws.Range("F15").Value = 282258108.13105428
$ws.Range("F15").Font.Name = "Arial"
$ws.Rows.Item(15).AutoFit()

$ws.Range("F16").Value = 82105216.118441164
$ws.Range("F16").Font.Name = "Arial"
$ws.Rows.Item(16).AutoFit()

$ws.Range("F17").Value = 428761167.17649448
$ws.Range("F17").Font.Name = "Arial"
$ws.Rows.Item(17).AutoFit()

$ws.Range("F18").Value = 341969516.72053099
$ws.Range("F18").Font.Name = "Arial"
$ws.Rows.Item(18).AutoFit()

$ws.Range("F19").Value = 67669214.883245558
$ws.Range("F19").Font.Name = "Arial"
$ws.Rows.Item(19).AutoFit()

$ws.Range("F20").Value = 23171201.136430159
$ws.Range("F20").Font.Name = "Arial"
$ws.Rows.Item(20).AutoFit()

$ws.Range("F21").Value = 34261412.749802783
$ws.Range("F21").Font.Name = "Arial"
$ws.Rows.Item(21).AutoFit()

$ws.Range("F22").Value = 374741182.86960781
$ws.Range("F22").Font.Name = "Arial"
$ws.Rows.Item(22).AutoFit()

$ws.Range("F23").Value = 48380324.546188213
$ws.Range("F23").Font.Name = "Arial"
$ws.Rows.Item(23).AutoFit()

$ws.Range("F24").Value = 1624889226.9504061
$ws.Range("F24").Font.Name = "Arial"
$ws.Rows.Item(24).AutoFit()

$ws.Range("F25").Value = 384603366.55867982
$ws.Range("F25").Font.Name = "Arial"
$ws.Rows.Item(25).AutoFit()

$ws.Range("F26").Value = 634255754.61794639
$ws.Range("F26").Font.Name = "Arial"
$ws.Rows.Item(26).AutoFit()

$ws.Range("F27").Value = 358978748.83680642
$ws.Range("F27").Font.Name = "Arial"
$ws.Rows.Item(27).AutoFit()

$ws.Range("F28").Value = 14095331.6287286
$ws.Range("F28").Font.Name = "Arial"
$ws.Rows.Item(28).AutoFit()

$ws.Range("F29").Value = 9885001.0605341773
$ws.Range("F29").Font.Name = "Arial"
$ws.Rows.Item(29).AutoFit()

$ws.Range("F30").Value = 13339518.422678489
$ws.Range("F30").Font.Name = "Arial"
$ws.Rows.Item(30).AutoFit()

$ws.Range("F31").Value = 2072384.2807455731
$ws.Range("F31").Font.Name = "Arial"
$ws.Rows.Item(31).AutoFit()

$ws.Range("F32").Value = 265799909.30403939
$ws.Range("F32").Font.Name = "Arial"
$ws.Rows.Item(32).AutoFit()

$ws.Range("F33").Value = 116868040.2482132
$ws.Range("F33").Font.Name = "Arial"
$ws.Rows.Item(33).AutoFit()

$ws.Range("F34").Value = 73018917.938560098
$ws.Range("F34").Font.Name = "Arial"
$ws.Rows.Item(34).AutoFit()

$ws.Range("F35").Value = 372687684.84123462
$ws.Range("F35").Font.Name = "Arial"
$ws.Rows.Item(35).AutoFit()

$ws.Range("F36").Value = 107918644.1539456
$ws.Range("F36").Font.Name = "Arial"
$ws.Rows.Item(36).AutoFit()

$ws.Range("F37").Value = 105170903.47280011
$ws.Range("F37").Font.Name = "Arial"
$ws.Rows.Item(37).AutoFit()

$ws.Range("F38").Value = 528615755.16888589
$ws.Range("F38").Font.Name = "Arial"
$ws.Rows.Item(38).AutoFit()

$ws.Range("F39").Value = 108785996.03378139
$ws.Range("F39").Font.Name = "Arial"
$ws.Rows.Item(39).AutoFit()

$ws.Range("F40").Value = 13836481.88611367
$ws.Range("F40").Font.Name = "Arial"
$ws.Rows.Item(40).AutoFit()

$ws.Range("F41").Value = 38790991.565392204
$ws.Range("F41").Font.Name = "Arial"
$ws.Rows.Item(41).AutoFit()

$ws.Range("F42").Value = 367177909.36214811
$ws.Range("F42").Font.Name = "Arial"
$ws.Rows.Item(42).AutoFit()

$ws.Range("F43").Value = 2108904091.805475
$ws.Range("F43").Font.Name = "Arial"
$ws.Rows.Item(43).AutoFit()

$ws.Range("F44").Value = 1723851404.0443709
$ws.Range("F44").Font.Name = "Arial"
$ws.Rows.Item(44).AutoFit()

$ws.Range("F45").Value = 138937105.51078349
$ws.Range("F45").Font.Name = "Arial"
$ws.Rows.Item(45).AutoFit()

$ws.Range("F46").Value = 912527170.95464969
$ws.Range("F46").Font.Name = "Arial"
$ws.Rows.Item(46).AutoFit()

$ws.Range("F47").Value = 736665252.97333908
$ws.Range("F47").Font.Name = "Arial"
$ws.Rows.Item(47).AutoFit()

$ws.Range("F48").Value = 1598293183.760381
$ws.Range("F48").Font.Name = "Arial"
$ws.Rows.Item(48).AutoFit()

$ws.Range("F49").Value = 121880591.05360229
$ws.Range("F49").Font.Name = "Arial"
$ws.Rows.Item(49).AutoFit()
